$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 5606  # was 3740
$ws.Range("I8").Value = 7343.3335  # was 3740
$ws.Range("J8").Value = 3000  # was 0
$ws.Range("K8").Value = 22030.0005  # was 11220
$ws.Range("L8").Value = 9000  # was 0
$ws.Range("M8").Value = -21891.0005  # was -11081
$ws.Range("N8").Value = -9278  # was empty (new cell)
# Row 40
$ws.Range("H40").Value = 1847.0588  # was 2100
$ws.Range("I40").Value = 1537.5  # was 2160
$ws.Range("J40").Value = 2122.2222  # was 2057.1428
$ws.Range("K40").Value = 1537.5  # was 2160
$ws.Range("L40").Value = 2122.2222  # was 2057.1428
$ws.Range("M40").Value = -1362.5  # was -1985
$ws.Range("N40").Value = -2472.2222  # was -2407.1428
# Row 43
$ws.Range("H43").Value = 846.5833  # was 743.93335
$ws.Range("I43").Value = 767.3333  # was 651
$ws.Range("J43").Value = 873  # was 758.2308
$ws.Range("K43").Value = 767.3333  # was 651
$ws.Range("L43").Value = 873  # was 758.2308
$ws.Range("M43").Value = -698.3333  # was -582
$ws.Range("N43").Value = -1011  # was -896.2308
# Row 82
$ws.Range("H82").Value = 2721  # was 1914
$ws.Range("I82").Value = 2721  # was 1914
$ws.Range("K82").Value = 8163  # was 5742
$ws.Range("M82").Value = -7757  # was -5336
# Row 85
$ws.Range("H85").Value = 2721  # was 1914
$ws.Range("I85").Value = 2721  # was 1914
$ws.Range("K85").Value = 8163  # was 5742
$ws.Range("M85").Value = -6759  # was -4338
# Row 86
$ws.Range("H86").Value = 1548.25  # was 1584.7142
$ws.Range("I86").Value = 1412.2858  # was 1432.1666
$ws.Range("K86").Value = 1412.2858  # was 1432.1666
$ws.Range("M86").Value = -289.2858000000001  # was -309.1666
# Row 89
$ws.Range("H89").Value = 1548.25  # was 1584.7142
$ws.Range("I89").Value = 1412.2858  # was 1432.1666
$ws.Range("K89").Value = 7061.429  # was 7160.833000000001
$ws.Range("M89").Value = -1445.429  # was -1544.833000000001
# Row 128
$ws.Range("H128").Value = 67389.5  # was 70000
$ws.Range("J128").Value = 67389.5  # was 70000
$ws.Range("L128").Value = 67389.5  # was 70000
$ws.Range("N128").Value = -77349.5  # was -79960
# Row 138
$ws.Range("H138").Value = 6051176  # was 5983196
$ws.Range("I138").Value = 1243170  # was 1197158.9
$ws.Range("J138").Value = 8067437  # was 8067437.5
$ws.Range("K138").Value = 3729510  # was 3591476.7
$ws.Range("L138").Value = 24202311  # was 24202312.5
$ws.Range("M138").Value = -3724370  # was -3586336.7
$ws.Range("N138").Value = -24212591  # was -24212592.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4150.766  # was 19304.555
$ws.Range("I32").Value = 3075.2666  # was 2124.3022
$ws.Range("J32").Value = 6048.706  # was 56242.1
$ws.Range("K32").Value = 3075.2666  # was 2124.3022
$ws.Range("L32").Value = 6048.706  # was 56242.1
$ws.Range("M32").Value = -2788.2666  # was -1837.3022
$ws.Range("N32").Value = -6622.706  # was -56816.1
# Row 74
$ws.Range("H74").Value = 8522.223  # was 5889.5557
$ws.Range("I74").Value = 2182.1667  # was 1653.8948
$ws.Range("J74").Value = 21202.334  # was 15949.25
$ws.Range("K74").Value = 2182.1667  # was 1653.8948
$ws.Range("L74").Value = 21202.334  # was 15949.25
$ws.Range("M74").Value = -1308.1667  # was -779.8948
$ws.Range("N74").Value = -22950.334  # was -17697.25
# Row 77
$ws.Range("H77").Value = 8522.223  # was 5889.5557
$ws.Range("I77").Value = 2182.1667  # was 1653.8948
$ws.Range("J77").Value = 21202.334  # was 15949.25
$ws.Range("K77").Value = 10910.8335  # was 8269.474
$ws.Range("L77").Value = 106011.67  # was 79746.25
$ws.Range("M77").Value = -6542.833500000001  # was -3901.474
$ws.Range("N77").Value = -114747.67  # was -88482.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 597.6667  # was 734.4
$ws.Range("I107").Value = 546.5  # was 624
$ws.Range("J107").Value = 700  # was 900
$ws.Range("K107").Value = 546.5  # was 624
$ws.Range("L107").Value = 700  # was 900
$ws.Range("M107").Value = 1373.5  # was 1296
$ws.Range("N107").Value = -4540  # was -4740
# Row 123
$ws.Range("H123").Value = 22991.428  # was 29750
$ws.Range("J123").Value = 22490  # was 31000
$ws.Range("L123").Value = 22490  # was 31000
$ws.Range("N123").Value = -32290  # was -40800
# Row 134
$ws.Range("H134").Value = 2751.9119  # was 2821.7666
$ws.Range("I134").Value = 1727.7727  # was 1776.7142
$ws.Range("J134").Value = 4629.5  # was 5260.222
$ws.Range("K134").Value = 5183.3181  # was 5330.142599999999
$ws.Range("L134").Value = 13888.5  # was 15780.666
$ws.Range("M134").Value = -2648.3181  # was -2795.142599999999
$ws.Range("N134").Value = -18958.5  # was -20850.666

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1742.8529  # was 1666.3529
$ws.Range("I58").Value = 1180.2  # was 1157.8462
$ws.Range("J58").Value = 3305.7778  # was 3319
$ws.Range("K58").Value = 1180.2  # was 1157.8462
$ws.Range("L58").Value = 3305.7778  # was 3319
$ws.Range("M58").Value = -977.2  # was -954.8462
$ws.Range("N58").Value = -3711.7778  # was -3725
# Row 107
$ws.Range("H107").Value = 659.52  # was 703.9545000000001
$ws.Range("I107").Value = 352.1111  # was 371.25
$ws.Range("J107").Value = 832.4375  # was 894.0714
$ws.Range("K107").Value = 352.1111  # was 371.25
$ws.Range("L107").Value = 832.4375  # was 894.0714
$ws.Range("M107").Value = 1567.8889  # was 1548.75
$ws.Range("N107").Value = -4672.4375  # was -4734.0714
# Row 132
$ws.Range("H132").Value = 3102.0527  # was 2953.9092
$ws.Range("I132").Value = 2090.7273  # was 1876.2858
$ws.Range("J132").Value = 4492.625  # was 4839.75
$ws.Range("K132").Value = 6272.1819  # was 5628.857400000001
$ws.Range("L132").Value = 13477.875  # was 14519.25
$ws.Range("M132").Value = -3742.1819  # was -3098.857400000001
$ws.Range("N132").Value = -18537.875  # was -19579.25
# Row 136
$ws.Range("H136").Value = 1742.8529  # was 1666.3529
$ws.Range("I136").Value = 1180.2  # was 1157.8462
$ws.Range("J136").Value = 3305.7778  # was 3319
$ws.Range("K136").Value = 3540.6  # was 3473.5386
$ws.Range("L136").Value = 9917.3334  # was 9957
$ws.Range("M136").Value = -990.6000000000004  # was -923.5385999999999
$ws.Range("N136").Value = -15017.3334  # was -15057

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 1101.3334  # was 574.5
$ws.Range("I46").Value = 0  # was 47.666668
$ws.Range("K46").Value = 0  # was 143.000004
$ws.Range("M46").ClearContents()  # was -52.00000399999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 17350.334  # was 23530.5
$ws.Range("I57").Value = 15000  # was 25000
$ws.Range("J57").Value = 22051  # was 22061
$ws.Range("K57").Value = 15000  # was 25000
$ws.Range("L57").Value = 22051  # was 22061
$ws.Range("M57").Value = -14180  # was -24180
$ws.Range("N57").Value = -23691  # was -23701
# Row 80
$ws.Range("H80").Value = 2644.182  # was 2850
$ws.Range("I80").Value = 2475  # was 2700
$ws.Range("J80").Value = 2740.8572  # was 2940
$ws.Range("K80").Value = 2475  # was 2700
$ws.Range("L80").Value = 2740.8572  # was 2940
$ws.Range("M80").Value = -1477  # was -1702
$ws.Range("N80").Value = -4736.8572  # was -4936
# Row 83
$ws.Range("H83").Value = 2644.182  # was 2850
$ws.Range("I83").Value = 2475  # was 2700
$ws.Range("J83").Value = 2740.8572  # was 2940
$ws.Range("K83").Value = 12375  # was 13500
$ws.Range("L83").Value = 13704.286  # was 14700
$ws.Range("M83").Value = -7383  # was -8508
$ws.Range("N83").Value = -23688.286  # was -24684
# Row 122
$ws.Range("H122").Value = 1389912.4  # was 1389913.6
$ws.Range("I122").Value = 1852734.9  # was 1588187
$ws.Range("J122").Value = 1445  # was 2000
$ws.Range("K122").Value = 5558204.699999999  # was 4764561
$ws.Range("L122").Value = 4335  # was 6000
$ws.Range("M122").Value = -5555754.699999999  # was -4762111
$ws.Range("N122").Value = -9235  # was -10900

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 7084.5  # was 7550.8667
$ws.Range("I22").Value = 451.1111  # was 463.44446
$ws.Range("J22").Value = 15613.143  # was 18182
$ws.Range("K22").Value = 451.1111  # was 463.44446
$ws.Range("L22").Value = 15613.143  # was 18182
$ws.Range("M22").Value = -156.1111  # was -168.44446
$ws.Range("N22").Value = -16203.143  # was -18772
# Row 27
$ws.Range("H27").Value = 7084.5  # was 7550.8667
$ws.Range("I27").Value = 451.1111  # was 463.44446
$ws.Range("J27").Value = 15613.143  # was 18182
$ws.Range("K27").Value = 451.1111  # was 463.44446
$ws.Range("L27").Value = 15613.143  # was 18182
$ws.Range("M27").Value = -344.1111  # was -356.44446
$ws.Range("N27").Value = -15827.143  # was -18396
# Row 40
$ws.Range("H40").Value = 4342.65  # was 4446
$ws.Range("I40").Value = 3093.25  # was 3331.3333
$ws.Range("K40").Value = 3093.25  # was 3331.3333
$ws.Range("M40").Value = -2957.25  # was -3195.3333
# Row 46
$ws.Range("H46").Value = 1341.8182  # was 1700
$ws.Range("I46").Value = 1112.5  # was 1900
$ws.Range("J46").Value = 1472.8572  # was 1500
$ws.Range("K46").Value = 1112.5  # was 1900
$ws.Range("L46").Value = 1472.8572  # was 1500
$ws.Range("M46").Value = -924.5  # was -1712
$ws.Range("N46").Value = -1848.8572  # was -1876
# Row 136
$ws.Range("H136").Value = 5778.6  # was 4879.391
$ws.Range("I136").Value = 2662.7273  # was 2036.5
$ws.Range("J136").Value = 9586.888999999999  # was 11377.429
$ws.Range("K136").Value = 7988.1819  # was 6109.5
$ws.Range("L136").Value = 28760.667  # was 34132.287
$ws.Range("M136").Value = -5438.1819  # was -3559.5
$ws.Range("N136").Value = -33860.667  # was -39232.287

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 677.1875  # was 714.0909
$ws.Range("I107").Value = 681.8889  # was 560.3333
$ws.Range("J107").Value = 671.1429000000001  # was 898.6
$ws.Range("K107").Value = 2045.6667  # was 1680.9999
$ws.Range("L107").Value = 2013.4287  # was 2695.8
$ws.Range("M107").Value = -125.6667000000002  # was 239.0001
$ws.Range("N107").Value = -5853.4287  # was -6535.8
# Row 128
$ws.Range("H128").Value = 49932.355  # was 46102.5
$ws.Range("J128").Value = 49932.355  # was 46102.5
$ws.Range("L128").Value = 49932.355  # was 46102.5
$ws.Range("N128").Value = -59892.355  # was -56062.5
